$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells (row 1): "_old" -> "_FV2310", "_new" -> "_FV2404" ---
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $txt = $cell.Text
    if ($txt -ne $null -and $txt -ne "") {
        $newTxt = $txt -replace '_old$', '_FV2310'
        $newTxt = $newTxt -replace '_new$', '_FV2404'
        if ($newTxt -ne $txt) {
            $cell.Value = $newTxt
        }
    }
}

# --- 2. Turn the used range into an Excel Table (ListObject) named "Table1" ---
$range = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
